$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe
$textUpdates = @(
    @('D2', '57.233.17'),
    @('E2', '  +3.65%  '),
    @('D3', '3.065.34'),
    @('E3', '  +6.12%  '),
    @('E4', '  +0.04%  '),
    @('E5', '  +4.89%  '),
    @('E6', '  +7.49%  '),
    @('E7', '  +0.02%  '),
    @('E8', '  +4.44%  '),
    @('E9', '  +1.23%  '),
    @('E10', '  +5.45%  '),
    @('E11', '  +7.70%  '),
    @('D12', '3.588.66'),
    @('E12', '  +6.20%  '),
    @('E13', '  +3.33%  '),
    @('E14', '  -0.01%  '),
    @('E15', '  +5.71%  '),
    @('D16', '57.254.90'),
    @('E16', '  +3.69%  '),
    @('D17', '3.068.85'),
    @('E17', '  +6.35%  '),
    @('E18', '  -0.05%  '),
    @('E19', '  +5.66%  '),
    @('E20', '  +8.05%  '),
    @('E21', '  +8.47%  '),
    @('E22', '  +0.35%  '),
    @('E23', '  +5.35%  '),
    @('E24', '  +5.83%  '),
    @('E25', '  +7.35%  '),
    @('D26', '0.0₃0950'),
    @('E26', '  +15.00%  '),
    @('E27', '  +0.75%  '),
    @('E28', '  +2.78%  '),
    @('E29', '  +2.31%  '),
    @('E30', '  +5.00%  '),
    @('E31', '  +6.16%  '),
    @('E32', '  +7.40%  '),
    @('E33', '  +4.06%  '),
    @('E34', '  +5.45%  '),
    @('E35', '  +6.32%  '),
    @('E36', '  +10.51%  '),
    @('E37', '  +5.62%  '),
    @('E38', '  +4.15%  '),
    @('D39', '3.104.65'),
    @('E39', '  +6.30%  '),
    @('E40', '  +2.55%  '),
    @('E41', '  +6.56%  '),
    @('E42', '  +0.06%  '),
    @('E43', '  +5.84%  '),
    @('D44', '2.236.75'),
    @('E44', '  +7.82%  '),
    @('E45', '  +10.47%  '),
    @('E46', '  +5.21%  '),
    @('E47', '  +5.04%  '),
    @('E49', '  +0.53%  '),
    @('E50', '  +3.99%  '),
    @('B51', 'SuiNetwork'),
    @('C51', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'),
    @('E51', '  +8.02%  ')
)

foreach ($u in $textUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# Numeric-looking values that must remain literal text (e.g. trailing zeros, 1.00 vs 1):
# force text format so Excel doesn't silently convert/round them to a Double
$numericTextUpdates = @(
    @('D5', '514.31'),
    @('D6', '140.62'),
    @('D7', '1.00'),
    @('D8', '0.434'),
    @('D9', '7.21'),
    @('D11', '0.371'),
    @('D14', '25.31'),
    @('D18', '5.93'),
    @('D19', '13.03'),
    @('D21', '336.55'),
    @('D23', '0.502'),
    @('D24', '65.32'),
    @('D25', '0.169'),
    @('D29', '7.03'),
    @('D30', '1.81'),
    @('D31', '20.70'),
    @('D33', '154.81'),
    @('D34', '4.54'),
    @('D36', '26.53'),
    @('D38', '0.0670'),
    @('D40', '36.94'),
    @('D41', '0.670'),
    @('D43', '3.83'),
    @('D45', '0.0252'),
    @('D46', '1.38'),
    @('D47', '0.946'),
    @('D48', '19.90'),
    @('D50', '0.0866'),
    @('D51', '0.687')
)

foreach ($u in $numericTextUpdates) {
    $r = $ws.Range($u[0])
    $r.NumberFormat = "@"
    $r.Value = $u[1]
    $r.Style = "Normal"
}

Write-Host "Applied $($textUpdates.Length + $numericTextUpdates.Length) cell updates"